$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix two existing cells where detect_structure (Q) flips from 1 to 0 ---
$ws.Cells.Item(60, 17).Value = 0
$ws.Cells.Item(70, 17).Value = 0

# --- Fix isPivot (O) on row 403 ---
$ws.Cells.Item(403, 15).Value = 2

# --- Backup (R) column: rows 405 and 406 go from blank/inline-string to numeric 0 ---
$ws.Cells.Item(405, 18).Value = 0
$ws.Cells.Item(406, 18).Value = 0

# --- Append new weekly rows 407-425 ---
$newRows = @(
    @(45474,602.9500122070312,644.8499755859375,601.2000122070312,643,643,11572602,2024,7,1,0,0,0,27,0,0,0),
    @(45481,641.5999755859375,673.7000122070312,628.5999755859375,651.1500244140625,651.1500244140625,11795322,2024,7,8,0,0,0,28,0,0,0),
    @(45488,652.5,662.75,637.5,643.4000244140625,643.4000244140625,5448576,2024,7,15,0,0,0,29,0,0,1),
    @(45495,640.0499877929688,731.4500122070312,611.25,724.9000244140625,724.9000244140625,20105706,2024,7,22,0,0,0,30,0,0,0),
    @(45502,724.9000244140625,742.2000122070312,698.0499877929688,728.5999755859375,728.5999755859375,8066844,2024,7,29,0,0,0,31,0,0,0),
    @(45509,705,746.6500244140625,697.3499755859375,739.25,739.25,8271132,2024,8,5,0,0,0,32,0,0,0),
    @(45516,731.0999755859375,741,702.6500244140625,720.0499877929688,720.0499877929688,5251512,2024,8,12,0,0,0,33,0,0,0),
    @(45523,720,752.4500122070312,714.5499877929688,731.9500122070312,731.9500122070312,8781826,2024,8,19,0,0,0,34,0,0,0),
    @(45530,732,754.75,715.4500122070312,753.1500244140625,753.1500244140625,9224388,2024,8,26,0,0,0,35,0,0,0),
    @(45537,755,772.9000244140625,746.0499877929688,750.6500244140625,750.6500244140625,5927135,2024,9,2,0,0,0,36,0,0,0),
    @(45544,749.9000244140625,766.7000122070312,740.3499755859375,755.25,755.25,5876886,2024,9,9,0,0,0,37,0,0,0),
    @(45551,751.2000122070312,770.9000244140625,742.8499755859375,768,768,5463815,2024,9,16,0,0,0,38,0,0,0),
    @(45558,768,795,761.5,789.5,789.5,7214011,2024,9,23,0,0,0,39,0,0,0),
    @(45565,793.7999877929688,796.7999877929688,748.9500122070312,755.75,755.75,3664716,2024,9,30,0,0,0,40,1,0,0),
    @(45572,755.75,764.4500122070312,726.2000122070312,742.5499877929688,742.5499877929688,5106629,2024,10,7,0,0,0,41,0,0,0),
    @(45579,742.3499755859375,759.4500122070312,723.1500244140625,745.3499755859375,745.3499755859375,4651183,2024,10,14,0,0,0,42,0,0,0),
    @(45586,748.4500122070312,776.8499755859375,701.2000122070312,742.8499755859375,742.8499755859375,6789669,2024,10,21,0,0,0,43,0,0,0),
    @(45593,740.5,773.9000244140625,732.4500122070312,744.9500122070312,744.9500122070312,3613171,2024,10,28,0,0,0,44,0,0,0),
    @(45600,744,744.7999877929688,702.2999877929688,710.3499755859375,710.3499755859375,4449906,2024,11,4,0,0,0,45,0,0,0)
)

$startRow = 407
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
    # Datetime column keeps the same date/time display format as the rest
    # of column A (reuses the workbook's existing "YYYY-MM-DD HH:MM:SS" style).
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    # Column R (18) stays blank/empty for the new rows (matches source data's
    # not-yet-backed-up marker used for freshly appended weeks).
}
